$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.952.94'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.892.32'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.37%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7341'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('E8').Value = '  -2.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.26'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06894'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7704'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07952'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.888.24'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.214'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.43'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.960.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.780'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '239.49'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007753'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.130.18'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.968'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.287'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.40'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1266'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.019'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -9.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.356'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.534'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.298'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.057'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05092'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.274'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7332'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.720'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01923'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.298'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '74.09'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4439'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.930'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8373'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.627'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.44%  '
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.753'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.037.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.35'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '935.58'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.91%  '
